$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New K (strikeout) values per row, replacing old Strike# values
$kValues = @{
    2 = 0
    4 = 1
    5 = 1
    6 = 0
    7 = 1
    8 = 0
    9 = 1
    10 = 2
    11 = 0
    12 = 1
    13 = 2
    14 = 1
    15 = 2
    16 = 2
    17 = 1
    18 = 1
    19 = 1
    20 = 2
    21 = 1
    22 = 0
    23 = 2
    24 = 1
    25 = 1
    26 = 0
    27 = 1
    28 = 1
    29 = 2
    30 = 1
    31 = 2
    32 = 0
    33 = 1
    34 = 2
    35 = 0
    36 = 2
    37 = 0
    38 = 1
    39 = 2
    40 = 0
    41 = 0
    42 = 0
    43 = 0
    44 = 1
    45 = 1
    46 = 1
    47 = 3
    48 = 0
    49 = 1
    50 = 2
    51 = 0
    52 = 0
    53 = 1
    54 = 0
    55 = 1
    56 = 2
    57 = 1
    58 = 2
    59 = 2
    60 = 0
    61 = 2
    62 = 1
    63 = 1
    64 = 0
    65 = 1
    66 = 1
    67 = 2
    68 = 0
    69 = 0
    70 = 2
    71 = 1
    72 = 0
    73 = 0
    74 = 2
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
